$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds a header row (row 1) followed by exposure-site records.
# Three new records need to become rows 2-4, pushing every existing record
# down by three rows (old row 2 -> new row 5, ..., old row 130 -> new row 133).
#
# We shift the data manually (copying cell-by-cell, bottom row first so we
# never clobber data before it's been read) rather than using
# Rows(...).Insert(), because Insert() copies the formatting of the row
# above into the freshly inserted rows - here that's the bold/centered
# header style - which would incorrectly style the new data rows.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = $lastRow; $r -ge 2; $r--) {
    $destRow = $r + 3
    $ws.Cells.Item($destRow, 1).Value = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($destRow, 2).Value = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($destRow, 3).Value = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($destRow, 4).Value = $ws.Cells.Item($r, 4).Value()
}

# Fill the now-vacant rows 2-4 with the new exposure site records.
$newRows = @(
    @("300 Grattan St, Parkville VIC 3050", -37.798908, 144.956176, "Melbourne (C)"),
    @("Bank St, Ascot Vale VIC 3032", -37.774352, 144.92733, "Moonee Valley (C)"),
    @("286 Buckley St, Essendon VIC 3040", -37.755176, 144.90226, "Moonee Valley (C)")
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = 2 + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
